$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title / date-range shared strings ---
$ws.Range("A8").Value = "Volume 30   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/24/2023  Through  7/30/2023"

# --- Crime statistics table updates (rows 14-30) ---

# Row 14
$ws.Range("M14").Value = 0

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = 200
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -33.333333333333

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = 25
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 76
$ws.Range("J16").Value = 69
$ws.Range("K16").Value = 10.144927536231
$ws.Range("L16").Value = 192.307692307692
$ws.Range("M16").Value = -37.190082644628
$ws.Range("N16").Value = -79.005524861878

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 129
$ws.Range("J17").Value = 138
$ws.Range("K17").Value = -6.521739130434
$ws.Range("L17").Value = 22.857142857142
$ws.Range("M17").Value = 41.758241758241
$ws.Range("N17").Value = -9.790209790209

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 58
$ws.Range("K18").Value = -44.827586206896
$ws.Range("L18").Value = -21.951219512195
$ws.Range("M18").Value = -77.777777777777
$ws.Range("N18").Value = -90.960451977401

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 13
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -43.478260869565
$ws.Range("I19").Value = 123
$ws.Range("J19").Value = 146
$ws.Range("K19").Value = -15.753424657534
$ws.Range("L19").Value = 61.842105263157
$ws.Range("M19").Value = -3.149606299212
$ws.Range("N19").Value = -45.814977973568

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 90
$ws.Range("J20").Value = 84
$ws.Range("K20").Value = 7.142857142857
$ws.Range("L20").Value = 66.666666666666
$ws.Range("M20").Value = 34.328358208955
$ws.Range("N20").Value = -92.531120331950

# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 17.647058823529
$ws.Range("F21").Value = 56
$ws.Range("H21").Value = -13.846153846153
$ws.Range("I21").Value = 464
$ws.Range("J21").Value = 511
$ws.Range("K21").Value = -9.197651663405
$ws.Range("L21").Value = 46.372239747634
$ws.Range("M21").Value = -17.730496453900
$ws.Range("N21").Value = -79.930795847750

# Row 22
$ws.Range("M22").Value = -25

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 11.111111111111
$ws.Range("I23").Value = 73
$ws.Range("J23").Value = 59
$ws.Range("K23").Value = 23.728813559322
$ws.Range("L23").Value = 32.727272727272
$ws.Range("M23").Value = 143.333333333333

# Row 24
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -35.294117647058
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 44
$ws.Range("H24").Value = -9.090909090909
$ws.Range("I24").Value = 340
$ws.Range("J24").Value = 311
$ws.Range("K24").Value = 9.324758842443
$ws.Range("L24").Value = 49.779735682819
$ws.Range("M24").Value = 7.594936708860

# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 83.333333333333
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = 6.25
$ws.Range("I25").Value = 222
$ws.Range("J25").Value = 228
$ws.Range("K25").Value = -2.631578947368
$ws.Range("L25").Value = 37.888198757764
$ws.Range("M25").Value = -22.105263157894

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = 200
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = -31.578947368421
$ws.Range("L26").Value = -13.333333333333

# Row 27
$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = -12.5
$ws.Range("L27").Value = -6.666666666666

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 12
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = -14.285714285714
$ws.Range("N28").Value = -25

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 12
$ws.Range("K29").Value = -20
$ws.Range("L29").Value = -14.285714285714
$ws.Range("M29").Value = 9.090909090909
$ws.Range("N29").Value = -20

# Row 30
$ws.Range("C30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("C30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("F30").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"
